# Populate the "Atlas_ResID" resource-id column (B) that was previously
# empty below the header row. Row 2 (EFT_INFO) points at the message-box
# icon atlas; the remaining effect rows (3-15) point at the shared
# "Ssetting" atlas.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = "msg_icon"

$ws.Range("B3:B15").Value2 = "Ssetting"

# Match the author's final selection/cursor position.
$ws.Range("E14").Select()
